$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "43.897.66"
$ws.Cells.Item(2, 5).Value = "  -0.92%  "

$ws.Cells.Item(3, 4).Value = "2.351.43"
$ws.Cells.Item(3, 5).Value = "  -0.92%  "

$ws.Cells.Item(4, 5).Value = "  -0.08%  "

$ws.Cells.Item(5, 4).Value = "0.674"
$ws.Cells.Item(5, 5).Value = "  -2.85%  "

$ws.Cells.Item(6, 4).Value = "240.76"
$ws.Cells.Item(6, 5).Value = "  -1.67%  "

$ws.Cells.Item(7, 4).Value = "73.09"
$ws.Cells.Item(7, 5).Value = "  -4.63%  "

$ws.Cells.Item(8, 5).Value = "  -0.04%  "

$ws.Cells.Item(9, 4).Value = "0.600"
$ws.Cells.Item(9, 5).Value = "  +0.66%  "

$ws.Cells.Item(10, 5).Value = "  -2.63%  "

$ws.Cells.Item(11, 4).Value = "58.93"
$ws.Cells.Item(11, 5).Value = "  +1.72%  "

$ws.Cells.Item(12, 4).Value = "33.29"
$ws.Cells.Item(12, 5).Value = "  +2.88%  "

$ws.Cells.Item(13, 5).Value = "  -0.08%  "

$ws.Cells.Item(14, 4).Value = "7.29"
$ws.Cells.Item(14, 5).Value = "  -2.82%  "

$ws.Cells.Item(15, 4).Value = "2.700.94"
$ws.Cells.Item(15, 5).Value = "  -0.90%  "

$ws.Cells.Item(16, 5).Value = "  -4.60%  "

$ws.Cells.Item(17, 4).Value = "0.908"
$ws.Cells.Item(17, 5).Value = "  -1.86%  "

$ws.Cells.Item(18, 4).Value = "2.356.13"
$ws.Cells.Item(18, 5).Value = "  -0.78%  "

$ws.Cells.Item(19, 4).Value = "43.805.60"
$ws.Cells.Item(19, 5).Value = "  -1.45%  "

$ws.Cells.Item(20, 4).Value = "0.0000103"
$ws.Cells.Item(20, 5).Value = "  -0.44%  "

$ws.Cells.Item(21, 4).Value = "6.72"
$ws.Cells.Item(21, 5).Value = "  +0.14%  "

$ws.Cells.Item(22, 4).Value = "78.05"
$ws.Cells.Item(22, 5).Value = "  -0.78%  "

$ws.Cells.Item(23, 4).Value = "256.12"
$ws.Cells.Item(23, 5).Value = "  -0.95%  "

$ws.Cells.Item(24, 4).Value = "1.94"
$ws.Cells.Item(24, 5).Value = "  +12.78%  "

$ws.Cells.Item(25, 5).Value = "  -0.06%  "

$ws.Cells.Item(26, 5).Value = "  +0.70%  "

$ws.Cells.Item(27, 5).Value = "  -3.13%  "

$ws.Cells.Item(28, 5).Value = "  -2.59%  "

$ws.Cells.Item(29, 4).Value = "2.27"
$ws.Cells.Item(29, 5).Value = "  -1.73%  "

$ws.Cells.Item(30, 5).Value = "  -2.36%  "

$ws.Cells.Item(31, 4).Value = "177.42"
$ws.Cells.Item(31, 5).Value = "  +1.06%  "

$ws.Cells.Item(32, 5).Value = "  -1.65%  "

$ws.Cells.Item(33, 5).Value = "  +0.64%  "

$ws.Cells.Item(34, 4).Value = "0.0752"
$ws.Cells.Item(34, 5).Value = "  -1.47%  "

$ws.Cells.Item(35, 4).Value = "5.47"
$ws.Cells.Item(35, 5).Value = "  +2.00%  "

$ws.Cells.Item(36, 4).Value = "5.13"
$ws.Cells.Item(36, 5).Value = "  -5.34%  "

$ws.Cells.Item(37, 5).Value = "  -3.08%  "

$ws.Cells.Item(39, 5).Value = "  -4.61%  "

$ws.Cells.Item(40, 4).Value = "0.0276"
$ws.Cells.Item(40, 5).Value = "  -0.42%  "

$ws.Cells.Item(41, 4).Value = "67.75"
$ws.Cells.Item(41, 5).Value = "  +26.09%  "

$ws.Cells.Item(42, 5).Value = "  +14.64%  "

$ws.Cells.Item(43, 5).Value = "  +8.16%  "

$ws.Cells.Item(44, 4).Value = "9.25"
$ws.Cells.Item(44, 5).Value = "  +1.18%  "

$ws.Cells.Item(45, 5).Value = "  -1.32%  "

$ws.Cells.Item(46, 5).Value = "  +3.77%  "

$ws.Cells.Item(47, 4).Value = "2.52"
$ws.Cells.Item(47, 5).Value = "  -0.32%  "

$ws.Cells.Item(48, 5).Value = "  -1.78%  "

$ws.Cells.Item(49, 5).Value = "  -0.01%  "

$ws.Cells.Item(50, 4).Value = "99.60"
$ws.Cells.Item(50, 5).Value = "  -2.76%  "

$ws.Cells.Item(51, 5).Value = "  -4.98%  "
